# DOMA-8525: ability to set Meter.isAutomatic via meters import
# Adds a new "Автоматический" column to the meter import example sheet,
# normalizes the Помещение/Лицевой счет/Номер счетчика columns to text,
# and fixes up a couple of stray values/date-typed cells so that every
# row in the sample follows the same text-based convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column T: "Автоматический" header, formatted like column S.
# ---------------------------------------------------------------------
$ws.Range("S1:S11").Copy() | Out-Null
$ws.Range("T1:T11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("T1").Value = "Автоматический"

# Match column T's width to the neighbouring R:S columns.
$ws.Range("T1").ColumnWidth = $ws.Range("S1").ColumnWidth

# ---------------------------------------------------------------------
# 2. Помещение (B), Лицевой счет (D), Номер счетчика (F) become text.
# ---------------------------------------------------------------------
$unitNumber = @{ 2 = "1"; 3 = "1"; 4 = "1"; 5 = "1"; 6 = "1"; 7 = "1"; 8 = "2"; 9 = "2"; 10 = "2"; 11 = "2" }
$account    = @{ 2 = "111"; 3 = "111"; 4 = "111"; 5 = "111"; 6 = "111"; 7 = "111"; 8 = "222"; 9 = "222"; 10 = "222"; 11 = "222" }
$meterNo    = @{ 2 = "1"; 3 = "1"; 4 = "2"; 5 = "2"; 6 = "2"; 7 = "2"; 8 = "11"; 9 = "11"; 10 = "22"; 11 = "33" }

foreach ($r in 2..11) {
    $bCell = $ws.Range("B$r")
    $bCell.NumberFormat = "@"
    $bCell.Value = $unitNumber[$r]

    $dCell = $ws.Range("D$r")
    $dCell.NumberFormat = "@"
    $dCell.Value = $account[$r]

    $fCell = $ws.Range("F$r")
    $fCell.NumberFormat = "@"
    $fCell.Value = $meterNo[$r]
}

# ---------------------------------------------------------------------
# 3. Показание 2/3/4 (I/J/K) blank cells: drop the stray "don't apply
#    number format" flag so they read as plain General cells.
# ---------------------------------------------------------------------
$generalBlanks = @("I2","J2","K2","I3","J3","K3","I4","J4","K4","I5","J5","K5","K6","K7","K8","I9","J9","K9","I10","J10","K10","I11","J11","K11")
foreach ($ref in $generalBlanks) {
    $ws.Range($ref).NumberFormat = "General"
}

# ---------------------------------------------------------------------
# 4. Дата поверки (M) switches from a real date value to the literal
#    text "2021-12-20" for every data row.
# ---------------------------------------------------------------------
foreach ($r in 2..11) {
    $mCell = $ws.Range("M$r")
    $mCell.NumberFormat = "@"
    $mCell.Value = "2021-12-20"
}

# Дата передачи показаний (L7) was a stray real date; normalise it to
# the same text value used throughout the rest of the column.
$l7 = $ws.Range("L7")
$l7.NumberFormat = "@"
$l7.Value = "2021-12-20"

# ---------------------------------------------------------------------
# 5. Показание 1 (H2) data fix: 100 -> 100.5
# ---------------------------------------------------------------------
$ws.Range("H2").Value = 100.5
